$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update Fecha (D3) and Volumen (J3)
$ws.Cells.Item(3, 4).Value = 44692
$ws.Cells.Item(3, 10).Value = 120

# Row 4: update Fecha, Volumen, Precio minimo/maximo/promedio, Unidad, Origen, Precio $/Kg, Kg o Unidades
$ws.Cells.Item(4, 4).Value = 44221
$ws.Cells.Item(4, 10).Value = 250
$ws.Cells.Item(4, 11).Value = 1300
$ws.Cells.Item(4, 12).Value = 1500
$ws.Cells.Item(4, 13).Value = 1420
$ws.Cells.Item(4, 14).Value = "$/atado"
$ws.Cells.Item(4, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(4, 16).Value = 1420
$ws.Cells.Item(4, 17).Value = 1

# Row 5: update Fecha, Volumen, Precio minimo/maximo/promedio, Unidad, Origen, Precio $/Kg, Kg o Unidades
$ws.Cells.Item(5, 4).Value = 44687
$ws.Cells.Item(5, 10).Value = 160
$ws.Cells.Item(5, 11).Value = 3000
$ws.Cells.Item(5, 12).Value = 3500
$ws.Cells.Item(5, 13).Value = 3250
$ws.Cells.Item(5, 14).Value = "$/docena de matas"
$ws.Cells.Item(5, 15).Value = "Región Metropolitana"
$ws.Cells.Item(5, 16).Value = 542
$ws.Cells.Item(5, 17).Value = 6
